$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L and M are formatted as Text (numFmtId 49, "@"), so assigning a
# numeric .Value directly would store it as a text string. Temporarily switch
# to a plain numeric format, assign the number, then restore the original
# Text format so the stored style index is unchanged but the value stays
# numeric.
function Set-NumericValue($cell, $value) {
    $fmt = $cell.NumberFormat()
    $cell.NumberFormat = "0"
    $cell.Value = $value
    $cell.NumberFormat = $fmt
}

# Row 609: update raw input C609; B609 (shared formula) recalculates automatically
$ws.Cells.Item(609, 3).Value = 56

# Row 610: update raw input C610; B610 recalculates automatically
$ws.Cells.Item(610, 3).Value = 42

# Row 611: update raw input C611; B611 recalculates automatically
$ws.Cells.Item(611, 3).Value = 60

# Row 612: fill in previously-empty data; B612/H612/J612/K612 are formulas
# that recalculate automatically once the raw inputs below are populated
$ws.Cells.Item(612, 3).Value = 66
$ws.Cells.Item(612, 5).Value = 2
$ws.Cells.Item(612, 6).Value = 2
$ws.Cells.Item(612, 7).Value = 8
Set-NumericValue $ws.Cells.Item(612, 12) 1
Set-NumericValue $ws.Cells.Item(612, 13) 0

# Row 613
$ws.Cells.Item(613, 3).Value = 30
$ws.Cells.Item(613, 5).Value = 2
$ws.Cells.Item(613, 6).Value = 2
$ws.Cells.Item(613, 7).Value = 10
Set-NumericValue $ws.Cells.Item(613, 12) 0
Set-NumericValue $ws.Cells.Item(613, 13) 0

# Row 614
$ws.Cells.Item(614, 3).Value = 31
$ws.Cells.Item(614, 5).Value = 2
$ws.Cells.Item(614, 6).Value = 2
$ws.Cells.Item(614, 7).Value = 12
Set-NumericValue $ws.Cells.Item(614, 12) 0
Set-NumericValue $ws.Cells.Item(614, 13) 0

# Row 615
$ws.Cells.Item(615, 3).Value = 49
$ws.Cells.Item(615, 5).Value = 2
$ws.Cells.Item(615, 6).Value = 2
$ws.Cells.Item(615, 7).Value = 10
Set-NumericValue $ws.Cells.Item(615, 12) 0
Set-NumericValue $ws.Cells.Item(615, 13) 0

# Row 616
$ws.Cells.Item(616, 3).Value = 74
$ws.Cells.Item(616, 5).Value = 4
$ws.Cells.Item(616, 6).Value = 4
$ws.Cells.Item(616, 7).Value = 12
Set-NumericValue $ws.Cells.Item(616, 12) 0
Set-NumericValue $ws.Cells.Item(616, 13) 0

# Row 617
$ws.Cells.Item(617, 3).Value = 43
$ws.Cells.Item(617, 5).Value = 4
$ws.Cells.Item(617, 6).Value = 3
$ws.Cells.Item(617, 7).Value = 11
Set-NumericValue $ws.Cells.Item(617, 12) 0
Set-NumericValue $ws.Cells.Item(617, 13) 0

# Row 618
$ws.Cells.Item(618, 3).Value = 6
$ws.Cells.Item(618, 5).Value = 4
$ws.Cells.Item(618, 6).Value = 4
$ws.Cells.Item(618, 7).Value = 10
Set-NumericValue $ws.Cells.Item(618, 12) 0
Set-NumericValue $ws.Cells.Item(618, 13) 0
